# Rename the header row of the "input" sheet to lowercase column names.
# New unique strings must be introduced in this order so that they are
# appended to the shared string table in the same order as the target
# workbook: parameter, description, value_1, value_2, distribution, unit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "description"
$ws.Range("D1").Value = "value_1"
$ws.Range("E1").Value = "value_2"
$ws.Range("G1").Value = "distribution"
$ws.Range("C1").Value = "unit"
# F1 ("shift") and H1 ("site_specific") are unchanged.

# Move the active selection on the "input" sheet from I35 to C2.
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
